# The data rows (2-3 and 4-5) need to swap their Fecha/Calidad/Volumen/
# Precio.../Unidad/Precio $/Kg/Kg o Unidades values pairwise, i.e.
#   row 2 <-> row 3
#   row 4 <-> row 5
# Columns A, B, C, E, F, G, H, O, R stay identical (unchanged by the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

Swap-Rows 2 3
Swap-Rows 4 5
